# Circle Language Spec: Replace the word 'globality' with 'module' in some
# more places (+ move the _GoBack bookmark to reflect the cursor position
# at save time).
#
# Strategy notes:
#  - Plain `$range.Text = "..."` replaces the text but keeps it inside the
#    single run that used to hold the found text (no run split).
#  - Word only splits a run when something about the run actually differs
#    at the boundary. Adding (and instantly removing) a same-point bookmark
#    is a clean, side-effect-free way to force such a split because a
#    bookmark start/end is its own empty element that must sit between two
#    runs - it leaves no residue behind once removed.

$d = $word.ActiveDocument

$script:tmpBookmarkSeq = 0
function Split-RunAt([int]$pos) {
    $script:tmpBookmarkSeq = $script:tmpBookmarkSeq + 1
    $name = "ZZZtmpSplit$($script:tmpBookmarkSeq)"
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add($name, $bmRange) | Out-Null
    $d.Bookmarks($name).Delete() | Out-Null
}

function Find-First([string]$text, [int]$searchFrom) {
    $r = $d.Range($searchFrom, $d.Content.End)
    $found = $r.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { return $null }
    return $r
}

# ---------------------------------------------------------------------
# Hunk 1: drop the old _GoBack bookmark that sits after the "Modules"
# heading text.
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# ---------------------------------------------------------------------
# Hunk 2: wrap "publically" - no text change, only (unreachable via COM)
# proofErr markers in real Word. We still split the run at the word's
# boundaries so the surrounding structure lines up as closely as possible.
# ---------------------------------------------------------------------
$r = Find-First "are just publically accessible" 0
$wStart = $r.Start + ("are just ").Length
$wEnd = $wStart + ("publically").Length
Split-RunAt $wStart
Split-RunAt $wEnd

# ---------------------------------------------------------------------
# Hunk 3: "globalities refer to eachother" -> "modules refer to eachother"
# ---------------------------------------------------------------------
$r = Find-First "You don’t need to let globalities refer to eachother" 0
$pStart = $r.Start
$prefixLen = ("You don’t need to let ").Length
$wStart = $pStart + $prefixLen
$wEnd = $wStart + ("globalities").Length
$rWord = $d.Range($wStart, $wEnd)
$rWord.Text = "modules"
# run split before/after "modules"
Split-RunAt $wStart
Split-RunAt ($wStart + ("modules").Length)
# run split around "eachother" (further down the same sentence)
$r2 = Find-First "refer to eachother" 0
$eaStart = $r2.Start + ("refer to ").Length
$eaEnd = $eaStart + ("eachother").Length
Split-RunAt $eaStart
Split-RunAt $eaEnd

# ---------------------------------------------------------------------
# Hunk 4: italic "Globality" -> "Module", and the _GoBack bookmark moves
# into the middle of "expression" (after "exp").
# ---------------------------------------------------------------------
$r = Find-First "Globality" 0
$r.Text = "Module"

$r = Find-First "expression in a diagram" 0
$pos = $r.Start + ("exp").Length
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------
# Hunk 5: "Multiple Globality Levels" -> "Multiple Module Levels"
# ---------------------------------------------------------------------
$r = Find-First "Multiple Globality Levels" 0
$pStart = $r.Start
$wStart = $pStart + ("Multiple ").Length
$wEnd = $wStart + ("Globality").Length
$rWord = $d.Range($wStart, $wEnd)
$rWord.Text = "Module"
Split-RunAt $wStart
Split-RunAt ($wStart + ("Module").Length)

# ---------------------------------------------------------------------
# Hunk 6: "embedded globalities" -> "embedded modules"
# ---------------------------------------------------------------------
$r = Find-First "embedded globalities" 0
$wStart = $r.Start + ("embedded ").Length
$wEnd = $wStart + ("globalities").Length
$rWord = $d.Range($wStart, $wEnd)
$rWord.Text = "modules"
Split-RunAt $wStart
Split-RunAt ($wStart + ("modules").Length)

# ---------------------------------------------------------------------
# Hunk 7: "multiple globalities" -> "multiple modules"
# ---------------------------------------------------------------------
$r = Find-First "multiple globalities" 0
$wStart = $r.Start + ("multiple ").Length
$wEnd = $wStart + ("globalities").Length
$rWord = $d.Range($wStart, $wEnd)
$rWord.Text = "modules"
Split-RunAt $wStart
Split-RunAt ($wStart + ("modules").Length)

# ---------------------------------------------------------------------
# Hunk 8: split "eachother" / "containmentwise" off, and
# "Globality is" -> "Module is"
# ---------------------------------------------------------------------
$r = Find-First "reference eachother and can be" 0
$eaStart = $r.Start + ("reference ").Length
$eaEnd = $eaStart + ("eachother").Length
Split-RunAt $eaStart
Split-RunAt $eaEnd

$r = Find-First "containmentwise ancestors" 0
$cwStart = $r.Start
$cwEnd = $cwStart + ("containmentwise").Length
Split-RunAt $cwStart
Split-RunAt $cwEnd

$r = Find-First "Globality is" 0
$wStart = $r.Start
$wEnd = $wStart + ("Globality").Length
$rWord = $d.Range($wStart, $wEnd)
$rWord.Text = "Module"
Split-RunAt $wStart
Split-RunAt ($wStart + ("Module").Length)

# ---------------------------------------------------------------------
# Hunk 9: "Multiple Globality Levels" heading -> "Multiple Module Levels"
#         (second occurrence, the "System Procedures for Globalities"
#         heading + paragraph)
# ---------------------------------------------------------------------
$r = Find-First "System Procedures for Globalities" 0
$wStart = $r.Start + ("System Procedures for ").Length
$wEnd = $wStart + ("Globalities").Length
$rWord = $d.Range($wStart, $wEnd)
$rWord.Text = "Modules"
Split-RunAt $wStart

$r = Find-First "System procedures of globalities are the same as for other object symbols." 0
$pStart = $r.Start

$off = ("System procedures of ").Length
$wStart = $pStart + $off
$wEnd = $wStart + ("globalities").Length
$d.Range($wStart, $wEnd).Text = "modules"
Split-RunAt $wStart
Split-RunAt ($wStart + ("modules").Length)

$r = Find-First "The only restriction is that globalities can’t have lines to objects" 0
$off = ("The only restriction is that ").Length
$wStart = $r.Start + $off
$wEnd = $wStart + ("globalities").Length
$d.Range($wStart, $wEnd).Text = "modules"
Split-RunAt $wStart
Split-RunAt ($wStart + ("modules").Length)

$r = Find-First "objects can’t have lines to globalities" 0
$off = ("objects can’t have lines to ").Length
$wStart = $r.Start + $off
$wEnd = $wStart + ("globalities").Length
$d.Range($wStart, $wEnd).Text = "modules"
Split-RunAt $wStart
Split-RunAt ($wStart + ("modules").Length)

$r = Find-First "Globalities can have an object line" 0
$wStart = $r.Start
$wEnd = $wStart + ("Globalities").Length
$d.Range($wStart, $wEnd).Text = "Modules"
Split-RunAt $wStart
Split-RunAt ($wStart + ("Modules").Length)

Write-Output "done"
